$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A (Product Code) to remain Text even though values look numeric.
# Row 9 is left completely untouched by this edit, so it is excluded here.
$ws.Range("A2:A8").NumberFormat = "@"
$ws.Range("A10:A27").NumberFormat = "@"

# Rows 2-8 and 10 get a full refresh of product data (row 9 is left untouched).
$fullRows = @(
  ,@(2, "1162105", "Enamelled Birth Flower Brooch In A Gift Box", "studiohop", 2, "https://www.notonthehighstreet.com/studiohop/product/enamelled-birth-flower-brooch-in-a-gift-box", "https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1162105&displayFeedbackType=PRODUCT&timeFrame=ALL")
  ,@(3, "1202855", "Dog And Owner Personalised Walking Socks", "alphabetinteriors", 2, "https://www.notonthehighstreet.com/alphabetinteriors/product/dog-and-owner-personalised-walking-socks", "https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1202855&displayFeedbackType=PRODUCT&timeFrame=ALL")
  ,@(4, "1268629", "Humorous Slate Weather Reader - Novelty Home Decor - Hanging Decoration - Gift For Home", "cgbgiftware", 2, "https://www.notonthehighstreet.com/cgbgiftware/product/humorous-slate-weather-reader", "https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1268629&displayFeedbackType=PRODUCT&timeFrame=ALL")
  ,@(5, "1305530", "F1 Formula One Gift Drinks Coaster Set Of Five", "iconiccoasters", 2, "https://www.notonthehighstreet.com/iconiccoasters/product/f1-formula-one-gift-drinks-coaster-set", "https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1305530&displayFeedbackType=PRODUCT&timeFrame=ALL")
  ,@(6, "1308378", "Women's White Cotton Nightdress Sleeveless Pink Lizzie", "minilunn", 2, "https://www.notonthehighstreet.com/minilunn/product/women-s-white-cotton-nightdress-sleeveless-pink-lizzie", "https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1308378&displayFeedbackType=PRODUCT&timeFrame=ALL")
  ,@(7, "1317092", "Fuck Cancer Sterling Silver Morse Code Chain Bracelet", "charlieboots", 2, "https://www.notonthehighstreet.com/charlieboots/product/fuck-cancer-sterling-silver-morse-code-chain-bracelet", "https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1317092&displayFeedbackType=PRODUCT&timeFrame=ALL")
  ,@(8, "1347156", "70th Birthday 1955 Sixpence Coin Compact Mirror", "ellieellie", 2, "https://www.notonthehighstreet.com/ellieellie/product/70th-birthday-1953-sixpence-coin-compact-mirror", "https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1347156&displayFeedbackType=PRODUCT&timeFrame=ALL")
  ,@(10, "1364790", "Personalised Bookrest In Solid Oak", "mijmoj", 2, "https://www.notonthehighstreet.com/mijmoj/product/personalised-bookrest-in-solid-oak", "https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1364790&displayFeedbackType=PRODUCT&timeFrame=ALL")
)

foreach ($row in $fullRows) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = "'" + $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
}

# Rows 11-27: only the Product Code, Review Count and Feefo URL remain; the
# Product Title / Seller Slug / NOTHS URL columns are blanked out (kept as empty
# text cells, matching column C which has always been blank).
$emptyRows = @(
  ,@(11, "1373657", 2)
  ,@(12, "1381151", 2)
  ,@(13, "1385338", 2)
  ,@(14, "1420726", 3)
  ,@(15, "1431678", 2)
  ,@(16, "1476302", 2)
  ,@(17, "1489678", 2)
  ,@(18, "1489841", 3)
  ,@(19, "377170", 2)
  ,@(20, "462337", 2)
  ,@(21, "469358", 3)
  ,@(22, "748820", 2)
  ,@(23, "837767", 2)
  ,@(24, "868202", 2)
  ,@(25, "876141", 3)
  ,@(26, "905169", 4)
  ,@(27, "928794", 2)
)

foreach ($row in $emptyRows) {
  $r = $row[0]
  $sku = $row[1]
  $reviewCount = $row[2]
  $feefoUrl = "https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=" + $sku + "&displayFeedbackType=PRODUCT&timeFrame=ALL"
  $ws.Cells.Item($r, 1).Value = "'" + $sku
  $ws.Cells.Item($r, 2).Value = "'"
  $ws.Cells.Item($r, 4).Value = "'"
  $ws.Cells.Item($r, 5).Value = $reviewCount
  $ws.Cells.Item($r, 6).Value = "'"
  $ws.Cells.Item($r, 7).Value = $feefoUrl
}

# Row 27 is a brand new row, so (unlike rows 11-26, whose blank "Seller" cell in
# column C already existed beforehand) its Seller cell must be created here too.
$ws.Cells.Item(27, 3).Value = "'"

